$d = $word.ActiveDocument

# 1. "Tech Lead" line: merge ": Abi" + "gail Conway" -> ": Abigail Conway"
$d.Content.Find.Execute(": Abi" + "gail Conway", $true, $false, $false, $false, $false, $true, 1, $false, ": Abigail Conway", 2) | Out-Null

# 2. "Team Members" line: merge ": Piper" + " " + "Dehoyos " -> ": Piper Dehoyos "
$d.Content.Find.Execute(": Piper" + " " + "Dehoyos ", $true, $false, $false, $false, $false, $true, 1, $false, ": Piper Dehoyos ", 2) | Out-Null

# 3. Stand up meeting date/time line: merge fragmented runs into one
$d.Content.Find.Execute(" " + "2/" + "6" + " 1:30 " + "after" + " class" + " ", $true, $false, $false, $false, $false, $true, 1, $false, " 2/6 1:30 after class ", 2) | Out-Null

# 4. Piper's agenda bullet: replace placeholder with actual content
$d.Content.Find.Execute("Piper will finish tech documentation for " + "_______ (will be given after scrum meeting)", $true, $false, $false, $false, $false, $true, 1, $false, "Piper will finish tech documentation for overriding string method, migrations, and python shell for running queries.", 2) | Out-Null

# 5. Scott's agenda bullet: replace placeholder with actual content
$d.Content.Find.Execute("Scott will " + "finish tech documentation for" + " " + "_______" + " " + "(will be given after scrum meeting)", $true, $false, $false, $false, $false, $true, 1, $false, "Scott will finish tech documentation for one-to-one model and one-to-many model.", 2) | Out-Null

# 6. Abi's agenda bullet: replace placeholder with actual content
$d.Content.Find.Execute("Abi" + " will finish tech documentation for" + " " + "_______" + " " + "(will be given after scrum meeting)", $true, $false, $false, $false, $false, $true, 1, $false, "Abi will finish tech documentation for set up admin and create model. ", 2) | Out-Null

# 7. Remove trailing space in "Time for async..." bullet AND delete the
# trailing empty ListParagraph paragraph that follows it. Matching the
# paragraph mark itself (^p) in the Find text and folding it into the
# replacement merges the "Time for async..." paragraph with the following
# (empty) one while keeping the *first* paragraph's properties (its
# numPr/bullet), which is what the target document keeps; the body's very
# final paragraph mark can't be deleted outright, so this mark-consuming
# replace is how the empty trailing <w:p> actually disappears.
$d.Content.Find.Execute("Time for async / track prosses with documentation ^p", $true, $false, $false, $false, $false, $true, 1, $false, "Time for async / track prosses with documentation", 2) | Out-Null
